$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3484939204877833
$ws.Range("C2").Value = 0.06223709392628507
$ws.Range("E2").Value = 0.4171142783858812
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.00241649757100462
$ws.Range("I2").Value = 0.4435338466345726
$ws.Range("K2").Value = 0.3724996450426943
$ws.Range("O2").Value = 1.948616353775151
$ws.Range("B3").Value = 0.3058381716033693
$ws.Range("C3").Value = 0.05512112264007385
$ws.Range("E3").Value = 0.3639871998607731
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.002418869006671513
$ws.Range("I3").Value = 0.4513819614751267
$ws.Range("K3").Value = 0.3252874019270848
$ws.Range("O3").Value = 1.97563386966678
$ws.Range("B4").Value = 0.2795943953214532
$ws.Range("C4").Value = 0.05072847772716216
$ws.Range("E4").Value = 0.331447085941619
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.002420399871837093
$ws.Range("I4").Value = 0.456533088198583
$ws.Range("K4").Value = 0.2962156432450058
$ws.Range("O4").Value = 1.993797288093489
$ws.Range("B5").Value = 0.2688870992604961
$ws.Range("C5").Value = 0.04893259515850445
$ws.Range("E5").Value = 0.3182053690840547
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002421042578568836
$ws.Range("I5").Value = 0.4587156186331089
$ws.Range("K5").Value = 0.2843482909942736
$ws.Range("O5").Value = 2.001593960718282
$ws.Range("B6").Value = 0.2671084082717243
$ws.Range("C6").Value = 0.04863403853310899
$ws.Range("E6").Value = 0.3160076743161113
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.002421150440980777
$ws.Range("I6").Value = 0.4590830586501706
$ws.Range("K6").Value = 0.2823765122144266
$ws.Range("O6").Value = 2.002912415873908
$ws.Range("B7").Value = 0.2794500438148191
$ws.Range("C7").Value = 0.05070428143221761
$ws.Range("E7").Value = 0.3312684300033197
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.002420408463528723
$ws.Range("I7").Value = 0.4565621850905863
$ws.Range("K7").Value = 0.2960556777321983
$ws.Range("O7").Value = 1.993900839020199
$ws.Range("B8").Value = 0.333797544591846
$ws.Range("C8").Value = 0.05978838886584015
$ws.Range("E8").Value = 0.3987785087160489
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.002417299754088083
$ws.Range("I8").Value = 0.4461707963942061
$ws.Range("K8").Value = 0.3562384621875196
$ws.Range("O8").Value = 1.957604592225096
$ws.Range("B9").Value = 0.4399324716094384
$ws.Range("C9").Value = 0.07741555569937475
$ws.Range("E9").Value = 0.5318754523272133
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.002411794336553484
$ws.Range("I9").Value = 0.4284377070950871
$ws.Range("K9").Value = 0.4735774122494831
$ws.Range("O9").Value = 1.898966217579627
$ws.Range("B10").Value = 0.5176222422298054
$ws.Range("C10").Value = 0.09025234667555537
$ws.Range("E10").Value = 0.6302061334237266
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002408105828749936
$ws.Range("I10").Value = 0.4170310572084652
$ws.Range("K10").Value = 0.5593556118152776
$ws.Range("O10").Value = 1.863590576237073
$ws.Range("B11").Value = 0.5528993887349998
$ws.Range("C11").Value = 0.09606749079657106
$ws.Range("E11").Value = 0.6750822206415421
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.002406504388814363
$ws.Range("I11").Value = 0.412196141996148
$ws.Range("K11").Value = 0.5982818952272737
$ws.Range("O11").Value = 1.849184014603992
$ws.Range("B12").Value = 0.5662482166005418
$ws.Range("C12").Value = 0.09826600879426906
$ws.Range("E12").Value = 0.6920982966960167
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.002405908901362836
$ws.Range("I12").Value = 0.4104163791192654
$ws.Range("K12").Value = 0.6130082266179784
$ws.Range("O12").Value = 1.843972192905639
$ws.Range("B13").Value = 0.5633737559892324
$ws.Range("C13").Value = 0.09779267724718466
$ws.Range("E13").Value = 0.6884325545406256
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.002406036664469222
$ws.Range("I13").Value = 0.4107974062773749
$ws.Range("K13").Value = 0.6098372889792927
$ws.Range("O13").Value = 1.845083797492677
$ws.Range("B14").Value = 0.5539978059757686
$ws.Range("C14").Value = 0.09624843561852003
$ws.Range("E14").Value = 0.6764816836234075
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.002406455178808758
$ws.Range("I14").Value = 0.412048694239008
$ws.Range("K14").Value = 0.5994937275072516
$ws.Range("O14").Value = 1.848750345369155
$ws.Range("B15").Value = 0.5482534634613785
$ws.Range("C15").Value = 0.09530207926741241
$ws.Range("E15").Value = 0.6691644057418245
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002406712954271105
$ws.Range("I15").Value = 0.4128218069959004
$ws.Range("K15").Value = 0.5931561296164034
$ws.Range("O15").Value = 1.851027974972425
$ws.Range("B16").Value = 0.5153154509564786
$ws.Range("C16").Value = 0.08987181889470719
$ws.Range("E16").Value = 0.6272764316697561
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002408212021057565
$ws.Range("I16").Value = 0.4173541718051013
$ws.Range("K16").Value = 0.5568097319012395
$ws.Range("O16").Value = 1.864566106795422
$ws.Range("B17").Value = 0.4950921355914204
$ws.Range("C17").Value = 0.0865342538850058
$ws.Range("E17").Value = 0.6016178007943438
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.002409151199493031
$ws.Range("I17").Value = 0.4202254519424926
$ws.Range("K17").Value = 0.5344877255716654
$ws.Range("O17").Value = 1.873303979068027
$ws.Range("B18").Value = 0.4834542016864418
$ws.Range("C18").Value = 0.08461228299634627
$ws.Range("E18").Value = 0.5868731816924679
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002409698591981981
$ws.Range("I18").Value = 0.4219102466756688
$ws.Range("K18").Value = 0.5216398243412925
$ws.Range("O18").Value = 1.878488394740458
$ws.Range("B19").Value = 0.4795127824308736
$ws.Range("C19").Value = 0.08396114446814806
$ws.Range("E19").Value = 0.5818831842433951
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.002409885168378297
$ws.Range("I19").Value = 0.4224864040459657
$ws.Range("K19").Value = 0.5172882399758976
$ws.Range("O19").Value = 1.880270956319535
$ws.Range("B20").Value = 0.4972455693116444
$ws.Range("C20").Value = 0.08688978090134469
$ws.Range("E20").Value = 0.6043477920757852
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.002409050477230646
$ws.Range("I20").Value = 0.4199163502831205
$ws.Range("K20").Value = 0.5368648639778542
$ws.Range("O20").Value = 1.872357392204009
$ws.Range("B21").Value = 0.5567520222773226
$ws.Range("C21").Value = 0.09670211354446678
$ws.Range("E21").Value = 0.6799913210588642
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.002406331954364384
$ws.Range("I21").Value = 0.4116797717532741
$ws.Range("K21").Value = 0.6025322725529918
$ws.Range("O21").Value = 1.847666769842505
$ws.Range("B22").Value = 0.5955851659292932
$ws.Range("C22").Value = 0.1030943225424892
$ws.Range("E22").Value = 0.7295605129487797
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.002404619002945239
$ws.Range("I22").Value = 0.4065947224961768
$ws.Range("K22").Value = 0.6453666362063473
$ws.Range("O22").Value = 1.832950639184347
$ws.Range("B23").Value = 0.5748646727600715
$ws.Range("C23").Value = 0.09968459074235625
$ws.Range("E23").Value = 0.7030919027533002
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.002405527421128307
$ws.Range("I23").Value = 0.4092813709209615
$ws.Range("K23").Value = 0.6225129280305453
$ws.Range("O23").Value = 1.840674535449025
$ws.Range("B24").Value = 0.4962720370191676
$ws.Range("C24").Value = 0.08672905698011846
$ws.Range("E24").Value = 0.603113541773709
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.00240909599061524
$ws.Range("I24").Value = 0.4200559890680751
$ws.Range("K24").Value = 0.5357902053958128
$ws.Range("O24").Value = 1.872784843013861
$ws.Range("B25").Value = 0.4112693007616599
$ws.Range("C25").Value = 0.07266694058863266
$ws.Range("E25").Value = 0.4957817183553743
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.002413220849700692
$ws.Range("I25").Value = 0.4329508711142633
$ws.Range("K25").Value = 0.4419085807279259
$ws.Range("O25").Value = 1.913480708087192
